$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell C3 was pointing at the wrong shared string ("<examplereplacement>");
# correct it to point at "no example".
$ws.Range("C3").Value = "no example"
